{"js": "// Replace each two-digit multiplication equation in the document with its\n// updated version, per the authoritative old->new mapping below.\nconst replacements = [\n  [\"24\u00d757=1368\", \"83\u00d765=5395\"],\n  [\"13\u00d796=1248\", \"77\u00d795=7315\"],\n  [\"23\u00d738=874\", \"63\u00d785=5355\"],\n  [\"40\u00d759=2360\", \"86\u00d776=6536\"],\n  [\"43\u00d796=4128\", \"27\u00d733=891\"],\n  [\"33\u00d714=462\", \"92\u00d719=1748\"],\n  [\"87\u00d796=8352\", \"39\u00d764=2496\"],\n  [\"94\u00d756=5264\", \"62\u00d752=3224\"],\n  [\"89\u00d789=7921\", \"93\u00d728=2604\"],\n  [\"97\u00d794=9118\", \"74\u00d754=3996\"],\n  [\"26\u00d753=1378\", \"33\u00d711=363\"],\n  [\"47\u00d731=1457\", \"31\u00d724=744\"],\n  [\"55\u00d731=1705\", \"27\u00d715=405\"],\n  [\"19\u00d772=1368\", \"23\u00d749=1127\"],\n  [\"15\u00d777=1155\", \"15\u00d787=1305\"],\n  [\"71\u00d729=2059\", \"35\u00d773=2555\"],\n  [\"72\u00d732=2304\", \"21\u00d750=1050\"],\n  [\"62\u00d729=1798\", \"65\u00d764=4160\"],\n  [\"14\u00d790=1260\", \"60\u00d738=2280\"],\n  [\"48\u00d711=528\", \"56\u00d770=3920\"],\n  [\"74\u00d783=6142\", \"66\u00d739=2574\"],\n  [\"67\u00d791=6097\", \"12\u00d750=600\"],\n  [\"50\u00d749=2450\", \"67\u00d799=6633\"],\n  [\"44\u00d730=1320\", \"70\u00d760=4200\"],\n  [\"95\u00d716=1520\", \"90\u00d742=3780\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation in the document with its\n# updated version, per the authoritative old->new mapping below.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"24\u00d757=1368\", \"83\u00d765=5395\"),\n  @(\"13\u00d796=1248\", \"77\u00d795=7315\"),\n  @(\"23\u00d738=874\", \"63\u00d785=5355\"),\n  @(\"40\u00d759=2360\", \"86\u00d776=6536\"),\n  @(\"43\u00d796=4128\", \"27\u00d733=891\"),\n  @(\"33\u00d714=462\", \"92\u00d719=1748\"),\n  @(\"87\u00d796=8352\", \"39\u00d764=2496\"),\n  @(\"94\u00d756=5264\", \"62\u00d752=3224\"),\n  @(\"89\u00d789=7921\", \"93\u00d728=2604\"),\n  @(\"97\u00d794=9118\", \"74\u00d754=3996\"),\n  @(\"26\u00d753=1378\", \"33\u00d711=363\"),\n  @(\"47\u00d731=1457\", \"31\u00d724=744\"),\n  @(\"55\u00d731=1705\", \"27\u00d715=405\"),\n  @(\"19\u00d772=1368\", \"23\u00d749=1127\"),\n  @(\"15\u00d777=1155\", \"15\u00d787=1305\"),\n  @(\"71\u00d729=2059\", \"35\u00d773=2555\"),\n  @(\"72\u00d732=2304\", \"21\u00d750=1050\"),\n  @(\"62\u00d729=1798\", \"65\u00d764=4160\"),\n  @(\"14\u00d790=1260\", \"60\u00d738=2280\"),\n  @(\"48\u00d711=528\", \"56\u00d770=3920\"),\n  @(\"74\u00d783=6142\", \"66\u00d739=2574\"),\n  @(\"67\u00d791=6097\", \"12\u00d750=600\"),\n  @(\"50\u00d749=2450\", \"67\u00d799=6633\"),\n  @(\"44\u00d730=1320\", \"70\u00d760=4200\"),\n  @(\"95\u00d716=1520\", \"90\u00d742=3780\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
